$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new content is a purely numeric-looking string must keep their
# original "text" semantics (shared-string cell, not a numeric cell). Excel's
# Value setter auto-coerces digit-only strings to numbers unless the cell is
# pre-formatted as Text, so force Text format (same effect as typing an
# apostrophe prefix) before writing those values.
$textCells = @("C3", "C11", "C13", "C14", "C24", "C28")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C3").Value = "15431275"
$ws.Range("C11").Value = "15431025"
$ws.Range("C13").Value = "15431340"
$ws.Range("C14").Value = "15431027"
$ws.Range("C24").Value = "15431104"
$ws.Range("C28").Value = "160628670"

$ws.Range("E14").Value = "PASS"
$ws.Range("C27").Value = "RT00006665"

$msgChrome = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome=124.0.6367.62)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL30', ip: '10.100.111.30', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '19.0.2'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 124.0.6367.62, chrome: {chromedriverVersion: 124.0.6367.78 (a087f2dd364d..., userDataDir: C:\Users\PARTH~1.SHA\AppDat...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:57094}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 7ec2368bf4d09dbd5ac05b819908f030
*** Element info: {Using=id, value=lblServiceID}
'@

$ws.Range("F13").Value = $msgChrome

Write-Output "Applied RTE OCP result updates"
